$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before the old row 11 (the "verifyPrimeSubjects_Old" block),
# which pushes all subsequent rows down by 3 and makes room for a new
# "verifyPrimeClasses_Parent" test-data block covering the Parent role.
$ws.Rows("11:13").Insert()

$platforms = @("Web", "Android", "iOS")
for ($i = 0; $i -lt 3; $i++) {
    $r = 11 + $i
    $ws.Range("B$r").Value = $platforms[$i]
    $ws.Range("F$r").Value = "Class 5, Class 6, Class 7"
    $ws.Range("C$r").Value = "verifyPrimeClasses_Parent"
    $ws.Range("D$r").Value = "Parent"
    $ws.Range("E$r").Value = "Prime Classes"
}

# Update the view state to reflect the author's final selection/scroll position.
$excel.ActiveWindow.ScrollRow = 4
$ws.Range("C13").Select()
